$p = $ppt.ActivePresentation

# --- 1) Update the "Update automatically" date field text that appears on
#        the slide master and every slide layout's Date Placeholder from
#        11/23/13 -> 11/24/13. ---
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/23/13") {
                $tr.Text = "11/24/13"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# --- 2) Slide 1 title: "Space Invaders IOS :) " -> "Space Invaders " ---
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Space Invaders IOS :) ") {
        $shp.TextFrame.TextRange.Text = "Space Invaders "
    }
}
